# Add a "Save" column (H) to the s_vals sheet, matching the formatting of
# the existing "sum" header in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the header formatting (bold font, border, centered alignment)
# from G1 onto the new H1 header cell, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New "Save" data column values.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
